# Generate Report for Handback
#
# The localization-status report is refreshed once the de-de / zh-cn
# handback packages have come back "in sync" with en-US: the Status
# column moves from "Ready for handoff" to "Handed back: in sync with
# en-US", and the per-row "Latest Target File" / "Latest Handback File"
# / "Latest Handback DateTime" columns (I/J/K) get populated + linked.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

$statusText = "Handed back: in sync with en-US"

$mdUrlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e4d3906e55972f5f0550b4051683c7227a15b0c3/e2e/4936858d-86d6-47ef-8f8c-e20afacb894c.md"
$mdUrlB = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e4d3906e55972f5f0550b4051683c7227a15b0c3/e2e/4bb901d1-63df-47ea-bc6a-c52d91dbe8a8.md"
$mdNameA = "4936858d-86d6-47ef-8f8c-e20afacb894c.md"
$mdNameB = "4bb901d1-63df-47ea-bc6a-c52d91dbe8a8.md"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------

# Status column -> handed back
$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

# Latest Target File / Latest Handback File / Latest Handback DateTime
$zhcn.Range("J2").Value = "4936858d-86d6-47ef-8f8c-e20afacb894c.aecc512cbd3f042c78a30624b062381f260f7356.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-05 15:16:32"
$zhcn.Range("J3").Value = "4bb901d1-63df-47ea-bc6a-c52d91dbe8a8.891e228a5739af0cbe341d3464de4a565bb6c33b.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-09-05 15:16:32"

# Rebuild the hyperlinks collection so the two new "Latest Target File"
# links (I2/I3) sit alongside the existing source-file links (A2/A3),
# re-inserted in row order so relationship ids come out 2,3,4,5.
$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), $mdUrlA, [Type]::Missing, [Type]::Missing, $mdNameA)
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdUrlA, [Type]::Missing, [Type]::Missing, $mdNameA)
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), $mdUrlB, [Type]::Missing, [Type]::Missing, $mdNameB)
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $mdUrlB, [Type]::Missing, [Type]::Missing, $mdNameB)

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------

$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

$dede.Range("J2").Value = "4936858d-86d6-47ef-8f8c-e20afacb894c.aecc512cbd3f042c78a30624b062381f260f7356.de-de.xlf"
$dede.Range("K2").Value = "2016-09-05 15:16:41"
$dede.Range("J3").Value = "4bb901d1-63df-47ea-bc6a-c52d91dbe8a8.891e228a5739af0cbe341d3464de4a565bb6c33b.de-de.xlf"
$dede.Range("K3").Value = "2016-09-05 15:16:41"

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), $mdUrlA, [Type]::Missing, [Type]::Missing, $mdNameA)
$dede.Hyperlinks.Add($dede.Range("I2"), $mdUrlA, [Type]::Missing, [Type]::Missing, $mdNameA)
$dede.Hyperlinks.Add($dede.Range("A3"), $mdUrlB, [Type]::Missing, [Type]::Missing, $mdNameB)
$dede.Hyperlinks.Add($dede.Range("I3"), $mdUrlB, [Type]::Missing, [Type]::Missing, $mdNameB)

# ---------------------------------------------------------------------
# Column widths widened to fit the longer status text / new file names.
# ColumnWidth is quantized by the host to 1/6-character steps, so we
# feed it the character-width figure (XML width minus the standard
# 5/6 padding) that lands in the bucket closest to the recorded value.
# ---------------------------------------------------------------------

$wideStatus = 29.9777047293527 - (5/6)   # -> lands on 30 exactly, closest bucket
$wideFile   = 40 - (5/6)                  # -> lands on 40 exactly

$overview.Columns.Item(5).ColumnWidth = $wideStatus
$overview.Columns.Item(6).ColumnWidth = $wideStatus

$zhcn.Columns.Item(3).ColumnWidth = $wideStatus
$zhcn.Columns.Item(9).ColumnWidth = $wideFile
$zhcn.Columns.Item(10).ColumnWidth = $wideFile

$dede.Columns.Item(3).ColumnWidth = $wideStatus
$dede.Columns.Item(9).ColumnWidth = $wideFile
$dede.Columns.Item(10).ColumnWidth = $wideFile
